# Remove duplicated OM keys from test data
# Renumber the OM_Key values 1001 -> 1004 and 1002 -> 1005 across all sheets,
# and update the sheet selections so the workbook lands on the last sheet
# (Flag_Priority) as the active tab, mirroring the author's final click-through.

$wb = $excel.ActiveWorkbook

$wsExtract  = $wb.Worksheets.Item("WMT_Extract")
$wsCourt    = $wb.Worksheets.Item("Court_Reports")
$wsInst     = $wb.Worksheets.Item("Inst_Reports")
$wsWarr     = $wb.Worksheets.Item("Flag_Warr_4_n")
$wsUpw      = $wb.Worksheets.Item("Flag_Upw")
$wsODue     = $wb.Worksheets.Item("Flag_O_Due")
$wsPriority = $wb.Worksheets.Item("Flag_Priority")

# ---- WMT_Extract: row 2 / row 3 OM_Key (numeric) + OM_Key text column ----
$wsExtract.Range("K2").Value = 1004
$wsExtract.Range("K3").Value = 1005

# ---- Court_Reports: OM_Key numeric + composite OM_Key text ----
$wsCourt.Range("D2").Value = 1004
$wsCourt.Range("E2").Value = "1004|WMT|C"
$wsCourt.Range("D3").Value = 1005
$wsCourt.Range("E3").Value = "1005|WMT|Z"

# ---- Inst_Reports: OM_Key numeric + composite OM_Key text ----
$wsInst.Range("D2").Value = 1004
$wsInst.Range("E2").Value = "1004|WMT|C"
$wsInst.Range("D3").Value = 1005
$wsInst.Range("E3").Value = "1005|WMT|Z"

# ---- Flag_Warr_4_n: OM_Key text column (F) ----
$wsWarr.Range("F2").Value = "1004"
$wsWarr.Range("F3").Value = "1005"

# ---- Flag_Upw: OM_Key text column (F) ----
$wsUpw.Range("F2").Value = "1004"
$wsUpw.Range("F3").Value = "1005"

# ---- Flag_O_Due: OM_Key text column (F) ----
$wsODue.Range("F2").Value = "1004"
$wsODue.Range("F3").Value = "1005"

# ---- Flag_Priority: OM_Key text column (F) ----
$wsPriority.Range("F2").Value = "1004"
$wsPriority.Range("F3").Value = "1005"

# ---- Update the view/selection state on every sheet to match the final
#      click-through, ending with Flag_Priority as the active tab. ----
$wsExtract.Activate()
$wsExtract.Range("K4").Select()

$wsCourt.Activate()
$wsCourt.Range("E4").Select()

$wsInst.Activate()
$wsInst.Range("E4").Select()

$wsWarr.Activate()
$wsWarr.Range("F4").Select()

$wsUpw.Activate()
$wsUpw.Range("F4").Select()

$wsODue.Activate()
$wsODue.Range("F4").Select()

$wsPriority.Activate()
$wsPriority.Range("F4").Select()
